$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = '#### CXCL9 
C-X-C motif chemokine ligand 9, also known as monokine induced by gamma interferon (MIG) is one of the most age-associated chemokines [[1](https://doi.org/10.1038/s43587-021-00082-y), [2](https://doi.org/10.3389/fimmu.2023.1177611)]. It is also involved in cardiac aging and vascular function [[1](https://doi.org/10.1038/s43587-021-00082-y)], regulation of chronic inflammation and inflammaging [[1](https://doi.org/10.1038/s43587-021-00082-y), [2](https://doi.org/10.3389/fimmu.2023.1177611), [3](https://doi.org/10.1007/s00198-024-07160-y)]. Levels of this biomarker may be elevated in primary Sjögren''s syndrome [[4](https://doi.org/10.1186/s13075-023-03229-x)], Chagas disease [[5](https://doi.org/10.1186%2Fs40249-020-00663-w)], chronic renal failure [[6](https://doi.org/10.1007/s11357-022-00540-4)], during progression of sarcoidosis [[7](https://doi.org/10.1016%2Fj.rmed.2019.105822)], and also be associated with age-related changes in muscle strength [[3](https://doi.org/10.1007/s00198-024-07160-y)] and frailty [[8](https://doi.org/10.1159/000535109)] in the elderly. Its concentration also varies in different types of cancer and may be a marker of the effectiveness of cancer treatment [[9](https://doi.org/10.1016%2Fj.ctrv.2017.11.007)].'
$ws.Range("T2").WrapText = $true

$ws.Range("T3").Value = '#### CCL11
CC-motif chemokine ligand 11, also known as eosinophil chemotactic protein (Eotaxin-1) is a chemokine that acts as a chemoattractant for eosinophilic granulocytes [[1](https://doi.org/10.1073/pnas.95.11.6273)] and is involved in the development of numerous diseases, associated with eosinophils, such as atopic dermatitis, allergic rhinitis, sinusitis, asthma, ulcerative colitis [[2](https://doi.org/10.3389/fnagi.2017.00402), [3](https://doi.org/10.1007/s00405-012-1927-5), [4](https://doi.org/10.3389/fimmu.2023.1243537)], as well as in the response to parasitic infections [[5](https://doi.org/10.1111/tmi.12095)], vascular inflammation and atherosclerosis [[6](https://doi.org/10.1161/01.CIR.102.18.2185)]. Age- related increases in CCL11 are associated with cognitive impairment, elevated CCL11 levels have been observed in psychiatric disorders as well as increased risk of Alzheimer''s disease and multiple sclerosis [[7](https://doi.org/10.3390/ph13090230), [8](https://doi.org/10.3389/fpsyt.2018.00241)].'
$ws.Range("T3").WrapText = $true

$ws.Range("T4").Value = '#### IL27
Interleukin-27 is a multifunctional cytokine that combines pro- and anti-inflammatory properties [[1](https://doi.org/10.1038/nri1648)] and can both promote and inhibit the development of inflammatory diseases, cancer, and viral infections [[2](https://doi.org/10.3389/fimmu.2023.1191228), [3](https://doi.org/10.3389/fimmu.2021.678515)]. IL27 expression is abnormal in various inflammatory autoimmune diseases (systemic lupus erythematosus, rheumatoid arthritis, Sjögren''s syndrome, Behçet''s disease, inflammatory bowel disease, multiple sclerosis, systemic sclerosis, and type 1 diabetes) [[4](https://doi.org/10.3389/fimmu.2024.1366377), [5](https://doi.org/10.1080/08830185.2020.1840565)]. IL27 has neuroprotective activity in the retina and brain [[6](https://doi.org/10.4103/1673-5374.336134)]. In COVID-19, its levels are positively correlated with advanced age, prolonged hospitalization, and severe disease [[7](https://doi.org/10.1038/s41419-020-03151-z)].'
$ws.Range("T4").WrapText = $true

$ws.Range("T5").Value = '#### IL5
Interleukin-5 is a key cytokine that mediates eosinophil activation and is involved in the development of many eosinophil-related diseases such as asthma, allergic rhinitis, and eosinophilic esophagitis [[1](https://doi.org/10.1016/B978-0-7234-3691-1.00105-7), [2](https://doi.org/10.1016/j.immuni.2019.03.022), [3](https://doi.org/10.1093/intimm/dxp102)]. Suppression of IL5 action by humanized anti-IL-5 monoclonal antibody can be used to treat asthma and other respiratory diseases [[4](https://doi.org/10.2183/pjab.87.463), [5](https://doi.org/10.3389/fphys.2019.01514), [6](https://doi.org/10.1111/all.16303)]. Decreased IL5 levels may be associated with coronary heart disease, as well as with the risk of angina pectoris and acute myocardial infarction [[7](https://doi.org/10.1016/j.rec.2019.07.005)].'
$ws.Range("T5").VerticalAlignment = -4160
$ws.Range("T5").WrapText = $true

$ws.Range("T6").Value = '#### CSF1
Colony Stimulating Factor-1 is one of the best known proinflammatory cytokines responsible for a variety of inflammatory diseases [[1](https://doi.org/10.1016/j.biopha.2018.04.046)]. Its levels depend on age [[2](https://doi.org/10.1038/s41598-017-07698-4), [3](https://doi.org/10.3389/fimmu.2023.1177611)] and can be associated with bad habits, including smoking [[4](https://doi.org/10.1038/s41598-017-07698-4)], as well as with various diseases such as periodontitis [[4](https://doi.org/10.1038/s41598-017-07698-4)], chronic kidney disease [[5](https://doi.org/10.1007%2Fs11357-022-00540-4)], atherosclerosis and cancer [[6](https://doi.org/10.1016/j.coi.2005.11.006)].'
$ws.Range("T6").WrapText = $true

$ws.Range("T7").Value = '#### CCL2
C-C motif chemokine ligand 2, also known as monocyte chemotactic protein 1 (MCP-1) is a proinflammatory chemokine [[1](https://doi.org/10.1089%2Fjir.2008.0027)] associated with a wide range of diseases: rheumatoid arthritis [[2](https://europepmc.org/article/med/20017457)], diabetes [[3](https://doi.org/10.1016/j.cyto.2012.06.018)], congestive heart failure [[4](https://doi.org/10.1161/01.CIR.97.12.1136)], Schnitzler syndrome [[5](https://doi.org/10.1111/bjd.17334)], chronic traumatic encephalopathy [[6](https://doi.org/10.1186/s12974-020-02036-4)], dyslipidemia in cardiovascular and renal diseases [[7](https://doi.org/10.1016/j.cyto.2019.154858)] and in many other pathologies [[8](https://doi.org/10.3389/fimmu.2019.02759)]. Abnormal CCL2 levels in pregnancy may indicate risks of spontaneous abortion, preeclampsia and preterm labor [[9](https://doi.org/10.3389%2Ffimmu.2022.1053457)].'
$ws.Range("T7").WrapText = $true

$ws.Range("T8").Value = '#### IL1B
Interleukin-1 Beta is a cytokine that plays an important role in enhancing immune response, a defense mechanism against pathogens and short-term inflammation [[1](https://doi.org/10.2119/molmed.2014.00232), [2](https://doi.org/10.1016/j.obmed.2021.100375), [3](https://doi.org/10.1038/cmi.2016.25)]. IL1B levels are responsible for the progression of obesity and impaired insulin signaling [[4](https://doi.org/10.4161/21623945.2014.979661)] and are associated with metabolic syndrome, invasiveness and progression of cancer [[5](https://doi.org/10.1186/s41232-019-0101-5)]. Elevated IL1B concentration is associated with worsening reproductive dysfunction in women with polycystic ovary syndrome [[6](https://doi.org/10.1016/s1995-7645(13)60030-9)], as well as with the intensity of inflammation in inflammatory bowel disease [[7](https://doi.org/10.3389/fmed.2024.1307394)]. IL1B plays a key role in intervertebral disc degeneration [[8](https://doi.org/10.1016/j.jcm.2017.09.003)], and its elevated levels are associated with autoinflammatory diseases [[9](https://doi.org/10.1182/blood-2010-07-273417)].'
$ws.Range("T8").VerticalAlignment = -4160
$ws.Range("T8").WrapText = $true

$ws.Range("T9").Value = '#### IL6
Interleukin-6 is a multifunctional cytokine that plays an important role in the protection of the human body and in affecting chronic inflammation and autoimmunity [[1](https://doi.org/10.1002%2Fpro.5560060501), [2](https://doi.org/10.1101%2Fcshperspect.a016295)]. Its levels increase with age [[3](https://doi.org/10.1186%2Fs12979-016-0076-x)] and may be associated with the risk of a wide range of diseases, including cardiovascular diseases [[4](https://doi.org/10.1186/s12916-022-02446-6)], autoimmune diseases [[5](https://doi.org/10.3389/fimmu.2023.1255533)], cancer [[6](https://doi.org/10.1093%2Fgerona%2F61.6.575)], type 2 diabetes [[7](https://doi.org/10.1001/jama.286.3.327)], rheumatoid arthritis [[8](https://doi.org/10.1002/art.1780310614)], and mortality risk, particularly from cardiovascular diseases [[9](https://doi.org/10.1016/S0002-9343(99)00066-2)].'
$ws.Range("T9").WrapText = $true

$ws.Range("T10").Value = '#### GCSF
Granulocyte Colony-Stimulating Factor, also known as colony-stimulating factor 3 (CSF3), is a cytokine that plays a key role in the regulation of neutrophil production [[1](https://doi.org/10.1016/j.molimm.2015.12.006)] and stimulates the proliferation and differentiation of stem cells into mature granulocytes [[2](https://doi.org/10.1016/B978-1-4377-2662-6.00018-3)]. It affects oocyte maturation and embryonic development [[3](https://doi.org/10.1016/B978-0-323-90805-4.00017-1)]. Elevated GCSF levels may lead to neutrophil abnormalities in severe COVID-19 [[4](https://doi.org/10.1002/JLB.3CE0121-002R)]. GCSF is elevated in age-related macular degeneration [[5](https://iovs.arvojournals.org/article.aspx?articleid=2790739)], has a contrasting effect in arthritis and inflammatory bowel diseases [[6](https://doi.org/10.1002/iub.361)].'
$ws.Range("T10").VerticalAlignment = -4160
$ws.Range("T10").WrapText = $true

$ws.Range("T11").Value = '#### CXCL10
C-X-C motif chemokine ligand 10, also known as  interferon gamma-induced protein 10 (IP-10) is a chemokine actively involved in regulating immune responses [[1](https://doi.org/10.1016/j.autrev.2008.12.002)]. CXCL10 levels are increased in various autoimmune diseases such as rheumatoid arthritis, systemic lupus erythematosus, Sjögren''s syndrome, systemic sclerosis, idiopathic inflammatory myopathy, and Chagas disease [[1](https://doi.org/10.1016/j.autrev.2008.12.002), [2](https://doi.org/10.1016/j.cyto.2021.155640), [3](https://doi.org/10.1186/s40249-020-00663-w)]. CXCL10 is involved in response to infections including rhinovirus [[4](https://doi.org/10.1164/rccm.200911-1673OC)], coxsackie virus, hepatitis B and C viruses [[5](https://doi.org/10.1002/jmv.10431)]. CXCL10 protects against severe acute respiratory syndrome caused by coronavirus [[6](https://doi.org/10.1016/j.cytogfr.2011.06.001), [7](https://doi.org/10.1186/s12979-024-00430-7)]. The observed levels of CXCL10 were higher in elderly compared to young adults and were negatively associated with working memory performance and higher in people with Alzheimer''s disease [[8](https://doi.org/10.1016/j.neurobiolaging.2017.11.009)].'
$ws.Range("T11").VerticalAlignment = -4160
$ws.Range("T11").WrapText = $true

$ws.Range("T12").Value = '#### VEGFA
Vascular Endothelial Growth Factor A is a cytokine that actively participates in inflammatory processes and is involved in angiogenesis, lymphangiogenesis, and hematopoiesis [[1](https://doi.org/10.1016/j.cyto.2012.10.024)]. Elevated levels of circulating VEGF have been observed in ischemic heart disease, diabetes, reproductive and inflammatory diseases [[1](https://doi.org/10.1016/j.cyto.2012.10.024), [2](https://doi.org/10.1016/j.thromres.2008.06.014)]. Increased VEGF expression contributes to the pathogenesis of inflammatory lung diseases [[3](https://doi.org/10.1155/2015/387842)], chronic inflammation in age-related macular degeneration [[4](https://doi.org/10.1002/jcp.22708)], as well as the development of cataracts, neovascular and non-exudative ocular pathologies [[5](https://doi.org/10.15252/emmm.201505613)]. Inhibition of VEGFA in the umbilical cord has been associated with negative pregnancy outcomes [[6](https://doi.org/10.1016/j.placenta.2019.09.002)].'
$ws.Range("T12").VerticalAlignment = -4160
$ws.Range("T12").WrapText = $true

$ws.Range("T13").Value = '#### TNF
Tumor Necrosis Factor is a pro-inflammatory cytokine that exerts various effects in response to acute inflammation [[1](https://doi.org/10.1016/bs.acc.2022.06.002), [2](https://doi.org/10.2174/1389450120666190821161259)]. In the healthy brain, baseline TNF-α release modulates plasticity for memory and learning processes [[3](https://doi.org/10.1038/nature04671)]; in neuroinflammation, increased TNF-α expression is associated with cytotoxicity and memory impairment [[4](https://doi.org/10.1007/s11481-011-9287-2)]. TNF has been implicated in a wide range of diseases such as cancer, arthritis, diabetes and other inflammatory diseases [[5](https://doi.org/10.1016/j.preteyeres.2005.09.001), [6](https://doi.org/10.2174/1568010023344571)]. Abnormal levels of TNF have been observed in rheumatoid arthritis, inflammatory bowel disease, psoriatic arthritis, psoriasis and uveitis [[7](https://doi.org/10.3390/ijms22052719)]. TNF is involved in the pathogenesis of atherosclerosis, type 2 diabetes mellitus, and Alzheimer''s disease in the elderly. Age-related increases in systemic TNF levels are associated with frailty and mortality risks [[8](https://doi.org/10.1016/s0889-8561(02)00056-5)].'
$ws.Range("T13").VerticalAlignment = -4160
$ws.Range("T13").WrapText = $true

$ws.Range("T14").Value = '#### PDGFB
Platelet Derived Growth Factor B is a cytokine that plays an important role in the development and differentiation of the vascular wall [[1](https://doi.org/10.1016/j.jocn.2009.02.009)]. Its elevated levels are observed in systemic lupus erythematosus [[2](https://doi.org/10.1186/s13075-017-1238-8)] and can be associated with an increased risk of acute myocardial infarction [[3](https://doi.org/10.18632/aging.205413)]. Increased PDGFB levels were observed in patients with ankylosing spondylitis and were associated with disease progression [[4](https://doi.org/10.1002/jbmr.4751)].'
$ws.Range("T14").VerticalAlignment = -4160
$ws.Range("T14").WrapText = $true

$ws.Range("T15").Value = '#### IL8
IInterleukin-8, also known as C-X-C motif chemokine ligand 8 (CXCL8), is a proinflammatory chemokine that activates neutrophils [[1](https://pubmed.ncbi.nlm.nih.gov/8315568/)] and is clinically used as a biomarker of inflammation [[2](https://doi.org/10.1159/000514885)]. Spontaneous IL8 secretion is lower in the elderly than in the young [[3](https://doi.org/10.1016/0047-6374(94)90020-5)]; however, blood IL8 concentration increases with age [[4](https://doi.org/10.1186/s12979-019-0151-1)] and in response to pain, and may characterize age-related dysregulation of the immune system [[5](https://doi.org/10.1016/j.exger.2015.09.017)]. IL8 is involved in the establishment and maintenance of the inflammatory microenvironment of the damaged vascular wall in cadiovascular diseases [[6](https://doi.org/10.1093/cvr/cvp241)], and its increase is associated with disease severity and prognosis for COVID-19 patients [[7](https://doi.org/10.3389/fimmu.2020.602395), [8](https://doi.org/10.3390/v15020549)].'
$ws.Range("T15").VerticalAlignment = -4160
$ws.Range("T15").WrapText = $true

$ws.Range("T16").Value = '#### PDGFA
Platelet Derived Growth Factor A is a cytokine that plays an important role in bone metabolism [[1](https://doi.org/10.1006/cyto.2000.0707)]. Decreased serum and pancreatic PDGFA levels contribute to decreased pancreatic function in aging [[2](https://doi.org/10.1016/j.bbrc.2019.12.057)]. PDGFA may be associated with different types of cancer [[3](https://doi.org/10.1006/excr.1993.1195), [4](https://doi.org/10.1016/0165-4608(94)90195-3)].'
$ws.Range("T16").VerticalAlignment = -4160
$ws.Range("T16").WrapText = $true

$ws.Range("T17").Value = '#### IL12Bp40
Interleukin-12 Beta subunit p40 is a cytokine that affects susceptibility to mycobacteria [[1](https://doi.org/10.1016/S0891-5520(03)00051-5)], acts as a chemoattractant for macrophages and promotes migration of bacteria-stimulated dendritic cells [[2](https://doi.org/10.1016/j.it.2006.11.002)]. The observed levels of IL12p40 increase with age [[3](https://doi.org/10.1006/cyto.1999.0537)]. IL12p40 is important for pathologies associated with psoriasis, inflammatory bowel disease and tumor growth [[4](https://doi.org/10.1016/j.cytogfr.2014.07.017), [5](https://doi.org/10.1111/exd.14566)], and IL12p40 deficiency abrogates IL17 production and prevents autoimmune response [[6](https://doi.org/10.1016/j.intimp.2006.09.024)]. Increased expression of IL12B has been observed in the central nervous system of patients with multiple sclerosis [[7](https://doi.org/10.1016/j.cyto.2017.10.024)], as well as in patients with hepatitis B [[8](https://doi.org/10.1016/j.meegid.2020.104663)].'
$ws.Range("T17").VerticalAlignment = -4160
$ws.Range("T17").WrapText = $true

$ws.Range("T18").Value = '#### IL15
Interleukin-15 is a proinflammatory cytokine [[1](https://doi.org/10.1016/B978-0-12-814841-9.00050-6), [2](https://doi.org/10.1016/j.autrev.2006.08.015)] that plays an important role in the development of inflammatory and protective immune responses by modulating immune cells of the innate and adaptive immune system [[3](https://doi.org/10.1016/j.micinf.2011.10.006)]. Its production is associated with increased osteoclastogenesis and bone destruction in rheumatoid arthritis [[1](https://doi.org/10.1016/B978-0-12-814841-9.00050-6)], and it contributes to the development and progression of various hematologic malignancies [[4](https://doi.org/10.3389/fimmu.2023.1141208)]. Elevated levels of IL15 can be observed in obese patients and serve as a risk marker for atherosclerosis [[5](https://doi.org/10.3389/fmed.2021.634962)]. Disturbed IL15 expression has been reported in patients with a number of inflammatory autoimmune diseases, including rheumatoid arthritis, sarcoidosis, chronic hepatitis C, ulcerative colitis, and psoriasis [[6](https://doi.org/10.1186/ar1202), [7](https://doi.org/10.1016/j.cyto.2020.155258)].'
$ws.Range("T18").VerticalAlignment = -4160
$ws.Range("T18").WrapText = $true

$ws.Range("T19").Value = '#### CXCL1
C-X-C motif chemokine ligand 1, also known as Growth-Regulated Oncogene alpha (GROa) is a cytokine, a neutrophil chemoattractant that plays an important role in the activation of innate immunity [[1](https://doi.org/10.1016/j.brainresbull.2021.05.024)]. CXCL1 induces regeneration of dopaminergic nerves [[2](https://doi.org/10.1016/j.brainresbull.2019.05.020)] and is elevated in chronic obstructive pulmonary disease and in smokers [[3](https://doi.org/10.1016/j.tips.2006.08.001)]. Plasma CXCL1 levels were elevated in patients with high-functioning autism spectrum disorder [[4](https://doi.org/10.1371/journal.pone.0020470)]. CXCL1 is involved in diseases of the cardiovascular system (atherosclerosis, atrial fibrillation, chronic ischemic heart disease, hypertension, sepsis), respiratory system (asthma, chronic obstructive pulmonary disease, rhinosinusitis, COVID-19, influenza, tuberculosis), skin (wound healing, psoriasis, burns), musculoskeletal system (osteoporosis, rheumatoid arthritis) [[5](https://doi.org/10.3390/ijms24010205), [6](https://doi.org/10.3390/ijms23084205), [7](https://doi.org/10.1038/s41420-023-01524-9)].'
$ws.Range("T19").VerticalAlignment = -4160
$ws.Range("T19").WrapText = $true

$ws.Range("T20").Value = '#### CCL4
C-C motif ligand 4, also known as Macrophage Inflammatory Protein 1 Beta (MIP-1b) is a proinflammatory chemokine that is released by activated cells in response to inflammation and regulates antitumor immunity [[1](https://doi.org/10.1016/j.semcancer.2022.06.010)]. Elevated levels of CCL4 are observed in systemic sclerosis, some infectious and bacterial diseases [[2](https://doi.org/10.1016/B978-0-08-102723-3.00252-3)]. CCL4 has been found to be associated with the pathogenesis of cardiovascular diseases, diabetes mellitus, autoimmune disorders, and HIV infection [[3](https://doi.org/10.1186/s12933-016-0439-9), [4](https://doi.org/10.1016/j.dci.2022.104525)]. CCL4 plays a key role in adipose tissue monocyte recruitment and demonstrate increased production in obesity [[5](https://doi.org/10.3390/ijms20184658)]. Circulating CCL4 increases with age [[6](https://doi.org/10.1007/s10456-024-09922-y)].'
$ws.Range("T20").VerticalAlignment = -4160
$ws.Range("T20").WrapText = $true

$ws.Range("T21").Value = '#### IFNA2
Interferon Alpha-2 is a pro-inflammatory cytokine secreted by virus-infected cells and acts to suppress viral infection [[1](https://doi.org/10.1016/j.gene.2015.04.087)]. Elevated IFNA2 levels are associated with systemic lupus erythematosus [[2](https://doi.org/10.1136/annrheumdis-2024-eular.4199), [3](https://doi.org/10.1177/0961203312468964)], rheumatoid arthritis [[4](https://doi.org/10.1038/s41598-024-72564-z)] and Sjögren''s syndrome [[5](https://doi.org/10.1093/rheumatology/keab688)], and are also increased in patients with COVID-19 [[6](https://doi.org/10.3389/fimmu.2023.1250214)].'
$ws.Range("T21").VerticalAlignment = -4160
$ws.Range("T21").WrapText = $true

$ws.Range("T22").Value = '#### IL13
Interleukin-13 is a cytokine with pro-inflammatory and anti-inflammatory properties, one of the central regulators of eosinophilic diseases [[1](https://doi.org/10.1016/j.anai.2012.06.014)], plays an important role in atopic disorders [[2](https://doi.org/10.1016/j.giec.2017.07.011), [3](https://doi.org/10.1111/all.13954)] and regulates many aspects of allergic inflammation [[3](https://doi.org/10.1111/all.13954), [4](https://doi.org/10.3389/fimmu.2018.00888)]. IL13 is involved in the pathogenesis of autoimmune diseases such as systemic lupus erythematosus, rheumatoid arthritis, systemic sclerosis, ulcerative colitis, type 1 diabetes, and Sjögren''s syndrome [[5](https://doi.org/10.1016/j.cytogfr.2018.12.001)]. IL13 levels are elevated in age-related macular degeneration [[6](https://doi.org/10.18240/ijo.2017.04.06)], and decrease with age in joints [[7](https://doi.org/10.1016/j.ocarto.2020.100128)].'
$ws.Range("T22").VerticalAlignment = -4160
$ws.Range("T22").WrapText = $true

$ws.Range("T23").Value = '#### FLT3L
Fms-related tyrosine kinase 3 ligand is a cytokine, an important regulator of hematopoiesis [[1](https://doi.org/10.3390/ijms18061115)]; its deficiency impairs monocyte development [[2](https://doi.org/10.1016/j.cell.2024.04.009)]. Elevated levels of FLT3L are observed in malaria and inflammatory bowel disease [[3](https://doi.org/10.1038/s41598-017-18283-0)], it is actively involved in acute myeloid leukemia [[4](https://doi.org/10.3390/ijms18061115), [5](https://doi.org/10.1016/j.autrev.2013.09.010)].'
$ws.Range("T23").VerticalAlignment = -4160
$ws.Range("T23").WrapText = $true

$ws.Range("T24").Value = '#### CD40LG
CD40 ligand is a cytokine that plays an important role in humoral responses, vascular physiopathology [[1](https://doi.org/10.1016/j.imbio.2011.03.010)], adaptive immunity [[2](https://doi.org/10.3389/fcvm.2017.00040)], and activation in response to inflammation [[3](https://doi.org/10.4049/jimmunol.168.2.713)]. It is involved in the effects of rheumatoid arthritis on bone [[4](https://doi.org/10.1016/B978-0-12-415853-5.00038-8)], atherosclerosis and other cardiovascular pathologies [[1](https://doi.org/10.1016/j.imbio.2011.03.010), [5](https://doi.org/10.3390/ijms21228533)]. CD40LG (and the CD40-CD40LG axis) is active in systemic lupus erythematosus [[6](https://doi.org/10.1016/j.clim.2009.05.011)], inflammatory bowel disease [[7](https://doi.org/10.3389/fimmu.2015.00529)], and other autoimmune diseases [[8](https://doi.org/10.1016/j.addr.2018.12.005)], as well as in central nervous system pathologies including brain injury, Alzheimer''s disease, Parkinson''s disease, stroke, epilepsy, and multiple sclerosis [[9](https://doi.org/10.3390/ijms23084115)].'
$ws.Range("T24").VerticalAlignment = -4160
$ws.Range("T24").WrapText = $true

$ws.Range("T25").Value = '#### CCL22
C-C motif ligand 22, also known as Macrophage-Derived Chemokine (MDC) is a chemokine involved in the transportation of cells to the focus of inflammation [[1](https://doi.org/10.1016/j.intimp.2020.106314)]. Decreased levels of CCL22 have been observed in patients with multiple sclerosis [[2](https://doi.org/10.1007/s10753-013-9775-z)], COVID-19 [[3](https://doi.org/10.3390/ijms241713083)]. Increased production of CCL22 is observed in parasitic infections [[4](https://doi.org/10.1016/j.cytogfr.2022.05.002)], psychotic disorders [[5](https://doi.org/10.1038/s41398-020-0776-z)].'
$ws.Range("T25").VerticalAlignment = -4160
$ws.Range("T25").WrapText = $true
